# C5-PowerPoint.pptx edit
#  1. Re-colour the deck's theme (the slide master's theme) from the
#     "Integral" palette to the stock "Office" palette - this is what
#     happens when a different Design/Theme swatch is applied from the
#     PowerPoint ribbon.
#  2. Re-style the table on slide 6 to use the built-in table style
#     {417FB4F1-6197-477A-9BD5-96FA473E8E51} instead of the custom
#     "Table_0" style.

function Set-ThemeColor($scheme, [int]$index, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $scheme.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation

# --- 1. Apply the "Office" colour scheme to the presentation theme -------
$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme

Set-ThemeColor $colorScheme 1  "000000"   # Dark 1
Set-ThemeColor $colorScheme 2  "FFFFFF"   # Light 1
Set-ThemeColor $colorScheme 3  "44546A"   # Dark 2
Set-ThemeColor $colorScheme 4  "E7E6E6"   # Light 2
Set-ThemeColor $colorScheme 5  "5B9BD5"   # Accent 1
Set-ThemeColor $colorScheme 6  "ED7D31"   # Accent 2
Set-ThemeColor $colorScheme 7  "A5A5A5"   # Accent 3
Set-ThemeColor $colorScheme 8  "FFC000"   # Accent 4
Set-ThemeColor $colorScheme 9  "4472C4"   # Accent 5
Set-ThemeColor $colorScheme 10 "70AD47"   # Accent 6
Set-ThemeColor $colorScheme 11 "0563C1"   # Hyperlink
Set-ThemeColor $colorScheme 12 "954F72"   # Followed hyperlink

# --- 2. Swap the "Sources of finance" table to a built-in table style ----
$slide = $p.Slides.Item(6)

for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shape = $slide.Shapes.Item($i)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle("{417FB4F1-6197-477A-9BD5-96FA473E8E51}")
    }
}
